$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(565).Insert()

$ws.Range("A565").Value = 5
$ws.Range("B565").Value = "Macroferia Regional de Talca"
$ws.Range("C565").Value = "Maule"
$ws.Range("D565").Value = 45142
$ws.Range("E565").Value = 7
$ws.Range("F565").Value = 100114014
$ws.Range("G565").Value = "Betarraga"
$ws.Range("H565").Value = "Sin especificar"
$ws.Range("I565").Value = "Primera"
$ws.Range("J565").Value = 5000
$ws.Range("K565").Value = 500
$ws.Range("L565").Value = 500
$ws.Range("M565").Value = 500
$ws.Range("N565").Value = "$/paquete 5 unidades"
$ws.Range("O565").Value = "Región del Maule"
$ws.Range("P565").Value = 100
$ws.Range("Q565").Value = 5
$ws.Range("R565").Value = "Hortaliza"
